# Fix locator for scenario 1
# Updates the "Expected Purchase Order" sheet: PO number, target-date window,
# sales-order-no locator, several forecast/inbound dates, inbound quantities,
# and the status of two rows (Completed -> Processing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / locator fields -------------------------------------------------

# Purchase Order No.
$ws.Cells.Item(14, 3).Value = "pDB204-2311001"

# Target Date (weekly window)
$ws.Cells.Item(19, 3).Value = "27 Nov 2023 - 03 Dec 2023"

# Forecast column header ("By Invoice Date" sub-range label)
$ws.Cells.Item(22, 14).Value = "04 Dec ~ 10 Dec"

# Sales Order No. locator for all three detail rows
$ws.Cells.Item(24, 5).Value = "sDB104-2311001"
$ws.Cells.Item(25, 5).Value = "sDB104-2311001"
$ws.Cells.Item(26, 5).Value = "sDB104-2311001"

# Status: Completed -> Processing
$ws.Cells.Item(25, 15).Value = "Processing"
$ws.Cells.Item(26, 15).Value = "Processing"

# --- Dates -------------------------------------------------------------

$ws.Cells.Item(18, 3).Value = 45250    # Order Date

$ws.Cells.Item(23, 16).Value = 45252   # P23 Inbounded Qty (Current Date) date
$ws.Cells.Item(23, 17).Value = 45301   # Q23 Inbound Plan Date
$ws.Cells.Item(23, 18).Value = 45343   # R23
$ws.Cells.Item(23, 19).Value = 45264   # S23
$ws.Cells.Item(23, 20).Value = 45266   # T23 Estimated Inbound Date

# --- Quantities ----------------------------------------------------------

$ws.Cells.Item(24, 16).Value = 0       # P24
$ws.Cells.Item(24, 20).Value = 800     # T24

$ws.Cells.Item(25, 16).Value = 0       # P25
$ws.Cells.Item(25, 20).Value = 1620    # T25 (was blank)

$ws.Cells.Item(26, 16).Value = 0       # P26
$ws.Cells.Item(26, 20).Value = 1620    # T26 (was blank)
